$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record was added at the top of the data (row 7), pushing the
# previously existing rows 7-17 down to rows 8-18.
$ws.Rows("7:7").Insert()

# Fill in the new row 7 with this week's data (mirrors the other rows'
# constant columns, with its own date/quality/volume/price figures).
$ws.Cells.Item(7, 1).Value = 1
$ws.Cells.Item(7, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(7, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(7, 4).Value = 44789
$ws.Cells.Item(7, 5).Value = 15
$ws.Cells.Item(7, 6).Value = "Fruta"
$ws.Cells.Item(7, 7).Value = 100101
$ws.Cells.Item(7, 8).Value = "Berries"
$ws.Cells.Item(7, 9).Value = 100101007
$ws.Cells.Item(7, 10).Value = "Kiwi"
$ws.Cells.Item(7, 11).Value = "Hayward"
$ws.Cells.Item(7, 12).Value = "Segunda"
$ws.Cells.Item(7, 13).Value = 250
$ws.Cells.Item(7, 14).Value = 19000
$ws.Cells.Item(7, 15).Value = 20000
$ws.Cells.Item(7, 16).Value = 19500
$ws.Cells.Item(7, 17).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(7, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(7, 19).Value = 1083
$ws.Cells.Item(7, 20).Value = 18
